# Apply "Updated symbol list" edit: refresh crypto price/volume snapshot values
# and fix the BKEXToken/CEJI/KickToken row ordering (rows 41-43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / Link cell updates (plain text, Excel will not reinterpret these) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

# --- Price / Volume(1h) cell updates ---
# These columns store plain numbers/percentages as TEXT in the workbook (t="inlineStr"),
# so we force the cells to Text format before assigning, otherwise Excel COM would
# auto-convert strings like "243.92" or "-0.70%" into real numbers/percentages.
# The style is reset back to "Normal" afterwards so no visible formatting changes.
$numericRefs = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "D20", "E21", "D22", "E22", "D23", "E23", "D25", "E25", "D26", "E26", "D27", "E27", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E47", "D48", "E48")
foreach ($ref in $numericRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "243.92"
$ws.Range("E2").Value = "-0.70%"
$ws.Range("D3").Value = "26.35"
$ws.Range("E3").Value = "3.29%"
$ws.Range("D4").Value = "5.136"
$ws.Range("E4").Value = "0.73%"
$ws.Range("D5").Value = "0.05594"
$ws.Range("E5").Value = "0.41%"
$ws.Range("D6").Value = "6.472"
$ws.Range("E6").Value = "-0.03%"
$ws.Range("D7").Value = "0.8202"
$ws.Range("E7").Value = "0.18%"
$ws.Range("D8").Value = "0.8344"
$ws.Range("E8").Value = "-1.15%"
$ws.Range("D9").Value = "0.1329"
$ws.Range("E9").Value = "-0.75%"
$ws.Range("D10").Value = "0.06994"
$ws.Range("E10").Value = "0.53%"
$ws.Range("D11").Value = "0.02888"
$ws.Range("E11").Value = "1.22%"
$ws.Range("D12").Value = "0.09385"
$ws.Range("E12").Value = "0.00%"
$ws.Range("D13").Value = "0.001514"
$ws.Range("E13").Value = "0.35%"
$ws.Range("D14").Value = "0.0006000"
$ws.Range("E14").Value = "-93.86%"
$ws.Range("D15").Value = "0.006137"
$ws.Range("E15").Value = "-0.09%"
$ws.Range("D16").Value = "3.656"
$ws.Range("E16").Value = "4.48%"
$ws.Range("E17").Value = "0.49%"
$ws.Range("D20").Value = "0.03108"
$ws.Range("E21").Value = "-2.22%"
$ws.Range("D22").Value = "3.736"
$ws.Range("E22").Value = "-0.80%"
$ws.Range("D23").Value = "0.04644"
$ws.Range("E23").Value = "-1.56%"
$ws.Range("D25").Value = "0.001246"
$ws.Range("E25").Value = "-0.30%"
$ws.Range("D26").Value = "0.004492"
$ws.Range("E26").Value = "-3.07%"
$ws.Range("D27").Value = "0.00009601"
$ws.Range("E27").Value = "-1.04%"
$ws.Range("E28").Value = "0.23%"
$ws.Range("D40").Value = "0.03643"
$ws.Range("E40").Value = "-0.45%"
$ws.Range("D41").Value = "0.006161"
$ws.Range("E41").Value = "-0.33%"
$ws.Range("D42").Value = "0.1052"
$ws.Range("E42").Value = "-0.02%"
$ws.Range("D43").Value = "0.002400"
$ws.Range("E43").Value = "-4.01%"
$ws.Range("D44").Value = "0.009001"
$ws.Range("E44").Value = "8.57%"
$ws.Range("D45").Value = "0.00005356"
$ws.Range("E45").Value = "0.96%"
$ws.Range("E47").Value = "8.24%"
$ws.Range("D48").Value = "0.002308"
$ws.Range("E48").Value = "8.71%"

foreach ($ref in $numericRefs) {
    $ws.Range($ref).Style = "Normal"
}
